# Apply updated cryptocurrency market data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the value to be stored as text even when it looks numeric
    # (e.g. "595.58"), without leaving a custom number-format style behind.
    if ($Text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $Range.NumberFormat = "@"
        $Range.Value = $Text
        $Range.Style = "Normal"
    } else {
        $Range.Value = $Text
    }
}

$ws.Range("D2").Value = '67.760.65'
$ws.Range("E2").Value = '  +2.07%  '

$ws.Range("D3").Value = '2.610.36'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("E4").Value = '  +0.27%  '

Set-TextValue $ws.Range("D5") '595.58'
$ws.Range("E5").Value = '  +0.66%  '

Set-TextValue $ws.Range("D6") '155.53'
$ws.Range("E6").Value = '  +0.77%  '

Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  +0.17%  '

Set-TextValue $ws.Range("D8") '0.548'
$ws.Range("E8").Value = '  +1.46%  '

$ws.Range("D9").Value = '2.606.59'
$ws.Range("E9").Value = '  +0.92%  '

Set-TextValue $ws.Range("D10") '0.126'
$ws.Range("E10").Value = '  +11.09%  '

$ws.Range("E11").Value = '  +0.86%  '

Set-TextValue $ws.Range("D12") '5.24'
$ws.Range("E12").Value = '  +0.27%  '

Set-TextValue $ws.Range("D13") '0.354'
$ws.Range("E13").Value = '  -0.60%  '

Set-TextValue $ws.Range("D14") '27.48'
$ws.Range("E14").Value = '  -2.99%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D15") '0.0000186'
$ws.Range("E15").Value = '  +2.85%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '3.087.27'
$ws.Range("E16").Value = '  +1.31%  '

$ws.Range("D17").Value = '67.653.67'
$ws.Range("E17").Value = '  +2.78%  '

$ws.Range("D18").Value = '2.615.13'
$ws.Range("E18").Value = '  +1.58%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D19") '11.18'
$ws.Range("E19").Value = '  -0.15%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D20") '364.88'
$ws.Range("E20").Value = '  +3.09%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '7.62'
$ws.Range("E21").Value = '  -2.88%  '

$ws.Range("E22").Value = '  -0.60%  '

Set-TextValue $ws.Range("D23") '2.01'
$ws.Range("E23").Value = '  -4.28%  '

Set-TextValue $ws.Range("D24") '1.00'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D25") '67.65'
$ws.Range("E25").Value = '  +2.17%  '

$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D26") '9.81'
$ws.Range("E26").Value = '  -6.07%  '

$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.742.06'
$ws.Range("E27").Value = '  +1.54%  '

Set-TextValue $ws.Range("D28") '0.0000103'
$ws.Range("E28").Value = '  -1.15%  '

$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D29") '576.61'
$ws.Range("E29").Value = '  -5.63%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.44%  '

$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D31") '1.42'
$ws.Range("E31").Value = '  -3.20%  '

Set-TextValue $ws.Range("D32") '7.93'
$ws.Range("E32").Value = '  -1.32%  '

Set-TextValue $ws.Range("D33") '1.86'
$ws.Range("E33").Value = '  +0.89%  '

Set-TextValue $ws.Range("D34") '0.133'
$ws.Range("E34").Value = '  -0.68%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D35") '0.998'
$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D36") '1.51'
$ws.Range("E36").Value = '  -4.47%  '

Set-TextValue $ws.Range("D37") '4.94'
$ws.Range("E37").Value = '  -1.87%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D38") '158.38'
$ws.Range("E38").Value = '  +2.13%  '

Set-TextValue $ws.Range("D39") '19.33'
$ws.Range("E39").Value = '  +1.03%  '

$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range("D40") '0.369'
$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D41") '1.85'
$ws.Range("E41").Value = '  +2.29%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D42") '5.33'
$ws.Range("E42").Value = '  -3.23%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D43") '2.54'
$ws.Range("E43").Value = '  -6.55%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D44") '41.20'
$ws.Range("E44").Value = '  -0.60%  '

$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D45") '16.43'
$ws.Range("E45").Value = '  +0.75%  '

$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D46") '0.999'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D47") '155.68'
$ws.Range("E47").Value = '  -0.54%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0289'
$ws.Range("E48").Value = '  -5.93%  '

Set-TextValue $ws.Range("D49") '3.73'
$ws.Range("E49").Value = '  -1.15%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D50") '0.631'
$ws.Range("E50").Value = '  +2.73%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D51") '20.77'
$ws.Range("E51").Value = '  -3.26%  '
